# Remove the trailing "Ver no Jupiter..." / copyright footer block that
# followed the "LOM3206: Eletrônica (Indicação de Conjunto)" requirement
# line: an empty paragraph, the "Ver no Jupiter Salvar em pdf Salvar em
# docx" paragraph, and the "© 2020 . Contact: ..." paragraph are deleted
# in their entirety, leaving the requirement paragraph directly followed
# by the (already-existing) empty paragraph that precedes the page break.

$d = $word.ActiveDocument

# Locate the "LOM3206" requirement paragraph via Find (avoids relying on
# a hard-coded paragraph index).
$found = $d.Content
$ok = $found.Find.Execute("LOM3206", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
$anchorStart = $found.Start

# Resolve which paragraph the match fell in.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($anchorStart -ge $p.Range.Start -and $anchorStart -lt $p.Range.End) {
        $anchorIndex = $i
        break
    }
}

# Delete the three paragraphs that immediately follow it (the blank
# paragraph, the "Ver no Jupiter..." line, and the "© 2020 ..." line).
$startPara = $d.Paragraphs.Item($anchorIndex + 1)
$endPara = $d.Paragraphs.Item($anchorIndex + 3)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
